$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2160.65
$ws.Range("J17").Value = 2160.65
$ws.Range("L17").Value = 6481.950000000001
$ws.Range("N17").Value = -6817.950000000001
$ws.Range("I33").Value = 35715176
$ws.Range("J33").Value = 1098
$ws.Range("K33").Value = 35715176
$ws.Range("L33").Value = 1098
$ws.Range("M33").Value = -35714947
$ws.Range("N33").Value = -1556
$ws.Range("H55").Value = 276.93332
$ws.Range("J55").Value = 376.4
$ws.Range("L55").Value = 376.4
$ws.Range("N55").Value = -804.4
$ws.Range("H94").Value = 551.4286
$ws.Range("I94").Value = 551.4286
$ws.Range("K94").Value = 551.4286
$ws.Range("M94").Value = -100.4286
$ws.Range("H125").Value = 18501.334
$ws.Range("J125").Value = 4494
$ws.Range("L125").Value = 40446
$ws.Range("N125").Value = -45366
$ws.Range("H132").Value = 2176.9412
$ws.Range("I132").Value = 1800.5333
$ws.Range("K132").Value = 5401.5999
$ws.Range("M132").Value = -2871.5999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5979.171
$ws.Range("I61").Value = 6101.359
$ws.Range("K61").Value = 6101.359
$ws.Range("M61").Value = -5889.359
$ws.Range("H136").Value = 5979.171
$ws.Range("I136").Value = 6101.359
$ws.Range("K136").Value = 18304.077
$ws.Range("M136").Value = -15754.077

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1006.125
$ws.Range("I22").Value = 1007
$ws.Range("K22").Value = 1007
$ws.Range("M22").Value = -834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 152.58333
$ws.Range("I7").Value = 80.09999999999999
$ws.Range("J7").Value = 515
$ws.Range("K7").Value = 80.09999999999999
$ws.Range("L7").Value = 515
$ws.Range("M7").Value = 32.90000000000001
$ws.Range("N7").Value = -741
$ws.Range("H41").Value = 18513.637
$ws.Range("J41").Value = 18513.637
$ws.Range("L41").Value = 18513.637
$ws.Range("N41").Value = -19369.637
$ws.Range("H51").Value = 14040
$ws.Range("J51").Value = 27600
$ws.Range("L51").Value = 27600
$ws.Range("N51").Value = -29072
$ws.Range("H58").Value = 13075.0625
$ws.Range("I58").Value = 9596
$ws.Range("K58").Value = 9596
$ws.Range("M58").Value = -9393
$ws.Range("H59").Value = 24459.4
$ws.Range("J59").Value = 20574.25
$ws.Range("L59").Value = 20574.25
$ws.Range("N59").Value = -22864.25
$ws.Range("H60").Value = 18562.3
$ws.Range("J60").Value = 20853.834
$ws.Range("L60").Value = 20853.834
$ws.Range("N60").Value = -21875.834
$ws.Range("H61").Value = 14040
$ws.Range("J61").Value = 27600
$ws.Range("L61").Value = 27600
$ws.Range("N61").Value = -28296
$ws.Range("H68").Value = 30025
$ws.Range("J68").Value = 26831.25
$ws.Range("L68").Value = 26831.25
$ws.Range("N68").Value = -28329.25
$ws.Range("H71").Value = 30025
$ws.Range("J71").Value = 26831.25
$ws.Range("L71").Value = 80493.75
$ws.Range("N71").Value = -87981.75
$ws.Range("H74").Value = 37523.332
$ws.Range("J74").Value = 37523.332
$ws.Range("L74").Value = 37523.332
$ws.Range("N74").Value = -39271.332
$ws.Range("H77").Value = 37523.332
$ws.Range("J77").Value = 37523.332
$ws.Range("L77").Value = 112569.996
$ws.Range("N77").Value = -121305.996
$ws.Range("H94").Value = 1759.8667
$ws.Range("I94").Value = 1801
$ws.Range("J94").Value = 1723.875
$ws.Range("K94").Value = 1801
$ws.Range("L94").Value = 1723.875
$ws.Range("M94").Value = -1350
$ws.Range("N94").Value = -2625.875
$ws.Range("H107").Value = 288.75
$ws.Range("J107").Value = 319.83334
$ws.Range("L107").Value = 319.83334
$ws.Range("N107").Value = -4159.83334
$ws.Range("H136").Value = 13075.0625
$ws.Range("I136").Value = 9596
$ws.Range("K136").Value = 28788
$ws.Range("M136").Value = -26238
$ws.Range("H141").Value = 182299.77
$ws.Range("J141").Value = 182299.77
$ws.Range("L141").Value = 182299.77
$ws.Range("N141").Value = -192659.77

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 750
$ws.Range("J25").Value = 468.33334
$ws.Range("L25").Value = 1405.00002
$ws.Range("N25").Value = -1743.00002
$ws.Range("H30").Value = 750
$ws.Range("J30").Value = 468.33334
$ws.Range("L30").Value = 1405.00002
$ws.Range("N30").Value = -1609.00002
$ws.Range("H50").Value = 534
$ws.Range("I50").Value = 394
$ws.Range("K50").Value = 1182
$ws.Range("M50").Value = -701
$ws.Range("H53").Value = 534
$ws.Range("I53").Value = 394
$ws.Range("K53").Value = 1182
$ws.Range("M53").Value = -701
$ws.Range("H60").Value = 111.166664
$ws.Range("I60").Value = 84.25
$ws.Range("K60").Value = 252.75
$ws.Range("M60").Value = -1.75
$ws.Range("H75").Value = 664.1818
$ws.Range("I75").Value = 763
$ws.Range("J75").Value = 642.2222
$ws.Range("K75").Value = 2289
$ws.Range("L75").Value = 1926.6666
$ws.Range("M75").Value = -1291
$ws.Range("N75").Value = -3922.6666
$ws.Range("H78").Value = 664.1818
$ws.Range("I78").Value = 763
$ws.Range("J78").Value = 642.2222
$ws.Range("K78").Value = 6867
$ws.Range("L78").Value = 5779.999800000001
$ws.Range("M78").Value = -1875
$ws.Range("N78").Value = -15763.9998
$ws.Range("H97").Value = 2230.1667
$ws.Range("I97").Value = 2095.75
$ws.Range("K97").Value = 6287.25
$ws.Range("M97").Value = -5791.25
$ws.Range("H109").Value = 10206.583
$ws.Range("I109").Value = 13636.143
$ws.Range("K109").Value = 40908.429
$ws.Range("M109").Value = -39868.429
$ws.Range("H112").Value = 200004910
$ws.Range("I112").Value = 250005140
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 750015420
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = -750014312
$ws.Range("N112").Value = -14216
$ws.Range("H123").Value = 4981.6665
$ws.Range("I123").Value = 4178
$ws.Range("J123").Value = 9000
$ws.Range("K123").Value = 12534
$ws.Range("L123").Value = 27000
$ws.Range("M123").Value = -10084
$ws.Range("N123").Value = -31900
$ws.Range("H131").Value = 14446844
$ws.Range("J131").Value = 2984.2
$ws.Range("L131").Value = 8952.599999999999
$ws.Range("N131").Value = -19032.6
$ws.Range("H134").Value = 1134.2106
$ws.Range("I134").Value = 1134.2106
$ws.Range("K134").Value = 3402.6318
$ws.Range("M134").Value = 1667.3682
$ws.Range("H139").Value = 8362.333000000001
$ws.Range("I139").Value = 3800.8823
$ws.Range("K139").Value = 11402.6469
$ws.Range("M139").Value = -6262.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 15000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -15312
$ws.Range("H122").Value = 2474.5
$ws.Range("I122").Value = 2459.2
$ws.Range("K122").Value = 7377.599999999999
$ws.Range("M122").Value = -4927.599999999999
$ws.Range("H126").Value = 7867.375
$ws.Range("I126").Value = 5316.6665
$ws.Range("K126").Value = 15949.9995
$ws.Range("M126").Value = -13479.9995
$ws.Range("H132").Value = 6261.2915
$ws.Range("I132").Value = 5430.222
$ws.Range("K132").Value = 16290.666
$ws.Range("M132").Value = -13760.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7997.1665
$ws.Range("I7").Value = 7997.1665
$ws.Range("K7").Value = 7997.1665
$ws.Range("M7").Value = -7885.1665
$ws.Range("H42").Value = 21331.615
$ws.Range("J42").Value = 39462.2
$ws.Range("L42").Value = 39462.2
$ws.Range("N42").Value = -40588.2
$ws.Range("H49").Value = 21331.615
$ws.Range("J49").Value = 39462.2
$ws.Range("L49").Value = 39462.2
$ws.Range("N49").Value = -39756.2
$ws.Range("H126").Value = 7997.1665
$ws.Range("I126").Value = 7997.1665
$ws.Range("K126").Value = 23991.4995
$ws.Range("M126").Value = -21521.4995
$ws.Range("H132").Value = 11287.4
$ws.Range("I132").Value = 11508.087
$ws.Range("J132").Value = 8749.5
$ws.Range("K132").Value = 34524.261
$ws.Range("L132").Value = 26248.5
$ws.Range("M132").Value = -31994.261
$ws.Range("N132").Value = -31308.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 32279.334
$ws.Range("J44").Value = 32279.334
$ws.Range("L44").Value = 32279.334
$ws.Range("N44").Value = -33387.334
$ws.Range("H52").Value = 31000
$ws.Range("J52").Value = 32000
$ws.Range("L52").Value = 32000
$ws.Range("N52").Value = -32452
$ws.Range("H122").Value = 3022.1538
$ws.Range("I122").Value = 1828.8
$ws.Range("K122").Value = 5486.4
$ws.Range("M122").Value = -3036.4
$ws.Range("H126").Value = 7396.8335
$ws.Range("I126").Value = 4811.5
$ws.Range("K126").Value = 14434.5
$ws.Range("M126").Value = -11964.5
$ws.Range("H136").Value = 8851
$ws.Range("I136").Value = 8031.2
$ws.Range("K136").Value = 24093.6
$ws.Range("M136").Value = -21543.6
